# Add SIN function description
# Reworks Sheet1 of the SIN function-description workbook:
#  - re-labels the existing "Number"/"SIN" header row, adding a third
#    "Formula  Text" column (D) that shows the formula text (via
#    FORMULATEXT) next to each existing SIN() example
#  - adds a second example table (rows 9-12) showing error propagation
#    through SIN() for a few invalid inputs, with a "Comments" column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): "Number" -> "angle", keep "SIN" in B1, add
#     "Formula  Text" (note: two spaces) in D1 --------------------------
$ws.Range("A1").Value = "angle"
$ws.Range("D1").Value = "Formula  Text"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats - reuse header style

# --- Column D: show the formula text used in column B for each row -----
$ws.Range("D2").Formula = '=FORMULATEXT($B2)'
# D3:D7 entered together as one range so Excel stores them as a shared
# formula group, same shape as the existing B4:B6 / B5 groups.
$ws.Range("D3:D7").Formula = '=FORMULATEXT($B3)'

# --- Second table header (row 9): Formula / Formula Text / (blank) / Comments
$ws.Range("A9").Value = "Formula"
$ws.Range("B9").Value = "Formula Text"
$ws.Range("D9").Value = "Comments"
$ws.Range("A1").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)   # xlPasteFormats - reuse header style

# --- Error-propagation examples (rows 10-12) ----------------------------
$ws.Range("A10").Formula = '=SIN(SQRT(-1))'
$ws.Range("B10").Formula = '=FORMULATEXT($A10)'
$ws.Range("D10").Value = "Example of error propagation."

$ws.Range("A11").Formula = '=SIN("str")'
$ws.Range("B11").Formula = '=FORMULATEXT($A11)'
$ws.Range("D11").Value = "Unable to convert angle argument to a number."

$ws.Range("A12").Formula = '=SIN(10/0)'
$ws.Range("B12").Formula = '=FORMULATEXT($A12)'
$ws.Range("D12").Value = "Input causes a #DIV/0! error."

# --- Selection ends on B13, as in the saved file ------------------------
$ws.Range("B13").Select() | Out-Null
